$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price / volume(1h) figures to match the latest scrape.
# Column D (Price) values are forced to text via a leading apostrophe so that
# Excel does not reinterpret numeric-looking strings as numbers (losing
# formatting / precision such as "1.001" or thousand-dot separated values).

$ws.Range('D2').Value = '''30.573.43'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').Value = '''1.877.62'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''238.99'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = '''0.4803'
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('D8').Value = '''0.2829'
$ws.Range('E8').Value = '  -1.92%  '
$ws.Range('D9').Value = '''0.06531'
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('D10').Value = '''1.869.00'
$ws.Range('E10').Value = '  -0.56%  '
$ws.Range('E11').Value = '  +1.34%  '
$ws.Range('D12').Value = '''16.60'
$ws.Range('E12').Value = '  -2.04%  '
$ws.Range('D13').Value = '''5.095'
$ws.Range('E13').Value = '  -2.21%  '
$ws.Range('D14').Value = '''88.63'
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('D15').Value = '''0.6622'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '''30.558.85'
$ws.Range('D17').Value = '''2.287.87'
$ws.Range('E17').Value = '  +6.90%  '
$ws.Range('D18').Value = '''13.33'
$ws.Range('E18').Value = '  -2.02%  '
$ws.Range('D19').Value = '''0.9998'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').Value = '''0.000007610'
$ws.Range('E20').Value = '  -1.60%  '
$ws.Range('D21').Value = '''226.37'
$ws.Range('E21').Value = '  +16.49%  '
$ws.Range('D22').Value = '''5.315'
$ws.Range('E22').Value = '  -2.81%  '
$ws.Range('D23').Value = '''1.001'
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = '''6.208'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').Value = '''9.317'
$ws.Range('E25').Value = '  -1.25%  '
$ws.Range('D26').Value = '''166.40'
$ws.Range('E26').Value = '  +0.97%  '
$ws.Range('D27').Value = '''18.67'
$ws.Range('E27').Value = '  +2.37%  '
$ws.Range('D28').Value = '''1.953'
$ws.Range('E28').Value = '  +1.29%  '
$ws.Range('D29').Value = '''1.464'
$ws.Range('E29').Value = '  +1.37%  '
$ws.Range('D30').Value = '''0.09496'
$ws.Range('E30').Value = '  +3.87%  '
$ws.Range('D31').Value = '''4.312'
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('D32').Value = '''4.031'
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('D33').Value = '''0.05020'
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('E34').Value = '  +6.42%  '
$ws.Range('E35').Value = '  +1.02%  '
$ws.Range('D36').Value = '''2.714'
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').Value = '''0.01829'
$ws.Range('E37').Value = '  -0.93%  '
$ws.Range('D38').Value = '''2.625'
$ws.Range('E38').Value = '  -0.20%  '
$ws.Range('D39').Value = '''2.072'
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('D40').Value = '''0.9073'
$ws.Range('E40').Value = '  -0.93%  '
$ws.Range('D41').Value = '''106.22'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''5.821'
$ws.Range('E42').Value = '  -1.05%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = '''0.4272'
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('D44').Value = '''1.005'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('D45').Value = '''7.475'
$ws.Range('E45').Value = '  -2.39%  '
$ws.Range('D46').Value = '''64.40'
$ws.Range('E46').Value = '  -1.32%  '
$ws.Range('D47').Value = '''0.1284'
$ws.Range('E47').Value = '  -4.51%  '
$ws.Range('D48').Value = '''1.473'
$ws.Range('E48').Value = '  -7.02%  '
$ws.Range('D49').Value = '''8.912'
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('D50').Value = '''33.77'
$ws.Range('E50').Value = '  -1.15%  '
$ws.Range('D51').Value = '''0.05649'
$ws.Range('E51').Value = '  -1.32%  '
